# Replace the data table content (rows 2-17) with updated exposure site data,
# and the header row text stays the same (Location / Site / Exposure period / Notes / Exist).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Location"
$ws.Range("B1").Value = "Site"
$ws.Range("C1").Value = "Exposure period"
$ws.Range("D1").Value = "Notes"
$ws.Range("E1").Value = "Exist"

$ws.Range("A2").Value = "Broadmeadows"
$ws.Range("B2").Value = "BonBon Bakery  Shop G134, Broadmeadows Central  1099/1168 Pascoe Vale Road  Broadmeadows VIC 3047"
$ws.Range("C2").Value = "12:30pm - 12:45pm 9/2/2021"
$ws.Range("D2").Value = "Case attended venue"
$ws.Range("E2").Value = "new"

$ws.Range("A3").Value = "Broadmeadows"
$ws.Range("B3").Value = "Craigieburn Line train"
$ws.Range("C3").Value = "1.25pm - 1.59pm  9/02/2021"
$ws.Range("D3").Value = "Case caught train from Broadmeadows Railway Station to Glenroy Railway Station"
$ws.Range("E3").Value = "old"

$ws.Range("A4").Value = "Broadmeadows"
$ws.Range("B4").Value = "Craigieburn Line train"
$ws.Range("C4").Value = "1:25pm - 1:59pm  9/02/2021"
$ws.Range("D4").Value = "Case caught train from Broadmeadows Railway Station to Glenroy Railway Station"
$ws.Range("E4").Value = "new"

$ws.Range("A5").Value = "Broadmeadows"
$ws.Range("B5").Value = "Ferguson Plarre Bakehouses - Broadmeadows, 1099-1169 Pascoe Vale Road"
$ws.Range("C5").Value = "12:30pm - 12:45pm 9/2/2021"
$ws.Range("D5").Value = "Case attended venue"
$ws.Range("E5").Value = "old"

$ws.Range("A6").Value = "Broadmeadows"
$ws.Range("B6").Value = "Woolworths  Broadmeadows Central  Pascoe Vale Road  Broadmeadows VIC 3047"
$ws.Range("C6").Value = "12:15pm - 12:30 pm 9/2/2021"
$ws.Range("D6").Value = "Case attended venue"
$ws.Range("E6").Value = "new"

$ws.Range("A7").Value = "Broadmeadows"
$ws.Range("B7").Value = "Woolworths Broadmeadows Central, Pascoe Vale Road"
$ws.Range("C7").Value = "12.15pm - 12:30 pm 9/2/2021"
$ws.Range("D7").Value = "Case attended venue"
$ws.Range("E7").Value = "old"

$ws.Range("A8").Value = "Glen Waverley"
$ws.Range("B8").Value = "Commonwealth Bank  28-32 Kingsway  Glen Waverley VIC 3015"
$ws.Range("C8").Value = "1:30pm-2:30pm 9/2/2021"
$ws.Range("D8").Value = "Case attended venue"
$ws.Range("E8").Value = "new"

$ws.Range("A9").Value = "Glen Waverley"
$ws.Range("B9").Value = "Commonwealth Bank, 28-32 Kingsway  Glen Waverley VIC 3015"
$ws.Range("C9").Value = "1:30pm-2:30pm 9/2/2021"
$ws.Range("D9").Value = "Case attended venue"
$ws.Range("E9").Value = "old"

$ws.Range("A10").Value = "Glen Waverley"
$ws.Range("B10").Value = "HSBC Bank,  38 Kingsway  Glen Waverley VIC 3015"
$ws.Range("C10").Value = "2:15pm-3:30pm 9/2/2021"
$ws.Range("D10").Value = "Case attended venue"
$ws.Range("E10").Value = "new"

$ws.Range("A11").Value = "Glen Waverley"
$ws.Range("B11").Value = "HSBC Bank, 38 Kingsway  Glen Waverley VIC 3015"
$ws.Range("C11").Value = "2:15pm-3:30pm 9/2/2021"
$ws.Range("D11").Value = "Case attended venue"
$ws.Range("E11").Value = "old"

$ws.Range("A12").Value = "Glenroy"
$ws.Range("B12").Value = "513 Eltham to Glenroy bus route  Glenroy Railway Station towards Eltham"
$ws.Range("C12").Value = "1.35pm  2.17pm  9/02/2021"
$ws.Range("D12").Value = "Case caught bus from Glenroy Railway Station towards Eltham"
$ws.Range("E12").Value = "old"

$ws.Range("A13").Value = "Glenroy"
$ws.Range("B13").Value = "513 Eltham to Glenroy bus route  Glenroy Railway Station towards Eltham"
$ws.Range("C13").Value = "1:35pm  2:17pm  9/02/2021"
$ws.Range("D13").Value = "Case caught bus from Glenroy Railway Station towards Eltham"
$ws.Range("E13").Value = "new"

$ws.Range("A14").Value = "Hoppers Crossing"
$ws.Range("B14").Value = "Coates Hire Werribee  148A Geelong Rd  Hoppers Crossing VIC 3029"
$ws.Range("C14").Value = "6.45am - 7.30am  8/02/21"
$ws.Range("D14").Value = "Case attended venue"
$ws.Range("E14").Value = "old"

$ws.Range("A15").Value = "Hoppers Crossing"
$ws.Range("B15").Value = "Coates Hire Werribee  148A Geelong Rd  Hoppers Crossing VIC 3029"
$ws.Range("C15").Value = "6:45am - 7:30am  8/02/21"
$ws.Range("D15").Value = "Case attended venue"
$ws.Range("E15").Value = "new"

$ws.Range("A16").Value = "Pascoe Vale"
$ws.Range("B16").Value = "Oak Park Sports and Aquatic Centre, 563a Pascoe Vale Road"
$ws.Range("C16").Value = "4pm - 7.30pm 10/2/2021"
$ws.Range("D16").Value = "Case attended venue"
$ws.Range("E16").Value = "old"

$ws.Range("A17").Value = "Pascoe Vale"
$ws.Range("B17").Value = "Oak Park Sports and Aquatic Centre, 563a Pascoe Vale Road"
$ws.Range("C17").Value = "4pm - 7:30pm 10/2/2021"
$ws.Range("D17").Value = "Case attended venue"
$ws.Range("E17").Value = "new"

# Resize columns to fit the new (longer) content, matching Excel's "best fit" widths
$ws.Columns.Item(1).ColumnWidth = 13.5
$ws.Columns.Item(2).ColumnWidth = 86.33333333333334
$ws.Columns.Item(3).ColumnWidth = 24.166666666666668
$ws.Columns.Item(4).ColumnWidth = 65.5

# Collapse the selection back down to the top-left cell
$ws.Range("A1").Select() | Out-Null
